$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the original B:AB values for every affected row BEFORE any writes,
# since several rows participate in multi-row cycles (not simple pairwise swaps).
$snap = @{}
$snap[20] = $ws.Range("B20:AB20").Value2
$snap[21] = $ws.Range("B21:AB21").Value2
$snap[38] = $ws.Range("B38:AB38").Value2
$snap[39] = $ws.Range("B39:AB39").Value2
$snap[48] = $ws.Range("B48:AB48").Value2
$snap[49] = $ws.Range("B49:AB49").Value2
$snap[81] = $ws.Range("B81:AB81").Value2
$snap[82] = $ws.Range("B82:AB82").Value2
$snap[130] = $ws.Range("B130:AB130").Value2
$snap[131] = $ws.Range("B131:AB131").Value2
$snap[132] = $ws.Range("B132:AB132").Value2
$snap[133] = $ws.Range("B133:AB133").Value2
$snap[135] = $ws.Range("B135:AB135").Value2
$snap[136] = $ws.Range("B136:AB136").Value2
$snap[157] = $ws.Range("B157:AB157").Value2
$snap[158] = $ws.Range("B158:AB158").Value2
$snap[159] = $ws.Range("B159:AB159").Value2
$snap[160] = $ws.Range("B160:AB160").Value2
$snap[162] = $ws.Range("B162:AB162").Value2
$snap[163] = $ws.Range("B163:AB163").Value2
$snap[164] = $ws.Range("B164:AB164").Value2

# Write each row its mapped source row's snapshot (column A/id stays untouched).
$ws.Range("B20:AB20").Value2 = $snap[21]
$ws.Range("B21:AB21").Value2 = $snap[20]
$ws.Range("B38:AB38").Value2 = $snap[39]
$ws.Range("B39:AB39").Value2 = $snap[38]
$ws.Range("B48:AB48").Value2 = $snap[49]
$ws.Range("B49:AB49").Value2 = $snap[48]
$ws.Range("B81:AB81").Value2 = $snap[82]
$ws.Range("B82:AB82").Value2 = $snap[81]
$ws.Range("B130:AB130").Value2 = $snap[131]
$ws.Range("B131:AB131").Value2 = $snap[133]
$ws.Range("B132:AB132").Value2 = $snap[130]
$ws.Range("B133:AB133").Value2 = $snap[132]
$ws.Range("B135:AB135").Value2 = $snap[136]
$ws.Range("B136:AB136").Value2 = $snap[135]
$ws.Range("B157:AB157").Value2 = $snap[159]
$ws.Range("B158:AB158").Value2 = $snap[157]
$ws.Range("B159:AB159").Value2 = $snap[158]
$ws.Range("B160:AB160").Value2 = $snap[162]
$ws.Range("B162:AB162").Value2 = $snap[160]
$ws.Range("B163:AB163").Value2 = $snap[164]
$ws.Range("B164:AB164").Value2 = $snap[163]
